$d = $word.ActiveDocument
$paras = $d.Paragraphs

# 1) Remove the red highlight ("Transferência de dados:", "POP", "PUSH").
#    Font.HighlightColorIndex applied on the full paragraph Range clears the
#    direct character formatting on both the run(s) and the paragraph mark.
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Transferência de dados:" -or $t -eq "POP" -or $t -eq "PUSH") {
        $p.Range.Font.HighlightColorIndex = 0
    }
}

# 2) Drop the stray en-US language mark on the "Aritimética:" paragraph
#    (paragraph mark + both runs) while keeping everything else (the
#    spell-check proofErr markers, rsids, run split) untouched.
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Aritimética:") {
        $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00566E8E" w:rsidRPr="00566E8E" w:rsidRDefault="00566E8E" w:rsidP="00566E8E"><w:pPr><w:ind w:firstLine="360"/></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00566E8E"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:t>Aritimética</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00566E8E"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/></w:rPr><w:t>:</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@
        $p.Range.InsertXML($xml)
    }
}
